$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "רפואה"
$ws.Range("A2").Value = "משכורת אריאנה (נטו)"
$ws.Range("A3").Value = "משכורת אלון (נטו)"
$ws.Range("A10").Value = "תחבורה"
$ws.Range("A21").Value = "אוכל בחוץ"
$ws.Range("A18").Value = "חדר כושר/חוגים"

$ws.Range("A18").Select()
